$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.476.92'
$ws.Range("E2").Value = '  +0.43%  '
$ws.Range("D3").Value = '3.508.54'
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.42'
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.02'
$ws.Range("E6").Value = '  +0.92%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +4.42%  '
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  +2.67%  '
$ws.Range("D12").Value = '4.104.87'
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000182'
$ws.Range("E14").Value = '  +1.41%  '
$ws.Range("D15").Value = '3.505.52'
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.90'
$ws.Range("E16").Value = '  -5.18%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '64.430.92'
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '9.96'
$ws.Range("E18").Value = '  +2.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.77'
$ws.Range("E19").Value = '  +3.36%  '
$ws.Range("E20").Value = '  -1.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '393.59'
$ws.Range("E21").Value = '  +2.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.576'
$ws.Range("E22").Value = '  +1.54%  '
$ws.Range("D23").Value = '3.648.23'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.65'
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("E28").Value = '  +0.10%  '
$ws.Range("E29").Value = '  -1.87%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.26'
$ws.Range("E30").Value = '  +2.24%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.29'
$ws.Range("E31").Value = '  +0.02%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.48'
$ws.Range("E32").Value = '  -5.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.158'
$ws.Range("E33").Value = '  +8.30%  '
$ws.Range("D34").Value = '3.531.74'
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.43'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.96'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("E39").Value = '  +0.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '168.18'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0789'
$ws.Range("E41").Value = '  +1.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.811'
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.25'
$ws.Range("E44").Value = '  -2.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.45'
$ws.Range("E45").Value = '  +1.64%  '
$ws.Range("E46").Value = '  +5.52%  '
$ws.Range("E47").Value = '  -2.97%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.80'
$ws.Range("E48").Value = '  +0.59%  '
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.398.89'
$ws.Range("E49").Value = '  -2.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.901'
$ws.Range("E50").Value = '  -0.63%  '
$ws.Range("E51").Value = '  +0.53%  '
